# FeedManagement.xlsx edit script
# Commit message: "Added cow birth status"
#
# Summary of the change:
#  - The single worksheet "Cow Feed management" is renamed to
#    "Cow Feed management-rose garden" and a near-duplicate worksheet
#    "Cow Feed management-kozhi panna" is appended (a second farm/location).
#  - A couple of small fixes/renames are applied before duplicating the
#    sheet so that both sheets inherit them:
#      * J3 formula corrected from SUM(I3:I9) to SUM(H3:H9)
#      * the "Total kgs" header relabelled "Total kgs/cow"
#  - After duplicating, "rose garden" gets its own tweaks (months-count
#    halved, "Kgs/month" header relabelled "Req Kgs/month" + wrapped)
#    while "kozhi panna" gets the left-hand data-entry table cleared and
#    a few input values/formulas adjusted to reflect its own herd.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Fixes applied on the original sheet BEFORE duplicating it, so both
#    resulting sheets carry them.
# ---------------------------------------------------------------------
$ws1.Range("J3").Formula = "=SUM(H3:H9)"

$ws1.Range("D2").Value = "Total kgs/cow"
$ws1.Range("J2").Value = "Total kgs/cow"
$ws1.Range("D10").Value = "Total kgs/cow"

# ---------------------------------------------------------------------
# 2) Duplicate the sheet (placed right after the original) so every
#    format/formula/merged-range is carried over faithfully.
# ---------------------------------------------------------------------
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 3) Rename both sheets.
# ---------------------------------------------------------------------
$ws1.Name = "Cow Feed management-rose garden"
$ws2.Name = "Cow Feed management-kozhi panna"

# ---------------------------------------------------------------------
# 4) "rose garden" (ws1) specific tweaks.
# ---------------------------------------------------------------------
$ws1.Range("M15").Value = 1

$ws1.Range("K16").Value = "Req Kgs/month"
$ws1.Range("K16").HorizontalAlignment = -4131
$ws1.Range("K16").WrapText = $true
$ws1.Rows.Item(16).RowHeight = 45

$ws1.Range("K17").Select()
$ws1.Application.ActiveWindow.ScrollRow = 2

# ---------------------------------------------------------------------
# 5) "kozhi panna" (ws2) specific tweaks.
# ---------------------------------------------------------------------
# Clear the left-hand "Milking cows" / "Dry cows" data-entry tables -
# this location tracks only the right-hand "Heifer" table for now.
$ws2.Range("A1:E7").ClearContents()
$ws2.Range("A9:E15").ClearContents()

$ws2.Range("J10").Value = 0
$ws2.Range("J11").Value = 0

$ws2.Range("H6").Value = 0.033
$ws2.Range("H7").Value = 0.033

$ws2.Range("K20").Formula = "=(H6*J12)*30"

$ws2.Range("H7").Select()

# ---------------------------------------------------------------------
# 6) Make "kozhi panna" the active tab, matching the saved view state.
# ---------------------------------------------------------------------
$ws2.Activate()
